$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276206556/TJaAJ2AJpb-camera-video-segment-driver-1743276209056.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOORQEFE7%2F20250331%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250331T040004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEDMaCXVzLXdlc3QtMiJIMEYCIQDYhldwn6xPibeo7UfO2hyaXchwZG%2BRUv3wHEvlAhzQ9wIhAKBCQOOXfTi9Kkihriivg46cyjPwSN9ea4juSW%2FfS77BKuYDCJz%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgwHR4J2QBwReOgfmJ8qugOPSbjSDvLSdXHFMzDd7rNCarOJd5s20MqO16y4tZIBDWKPK59sZrEyTMhwjPUPSLReD3kg2E1YPD%2Bvwc7OKmU7bmHhdFDUfIG7tKilS4GYE3JEvrsLY1pnH0X0VXlYEqm5ItI8QIHmvn29iLVQRQ1DWamIZa5%2BriHHXT9WhraXH2N4tTbrU8elcJX8DP%2FO8iLZCq5hMcBNdYb1O8WZTqxqCOZd1DjKZRtHgIQck5TdBZb0gwYLBPQzhamSl4kMMu%2FMH0vEeb%2B%2BZjz3r6f5V4I0Iy0gMIOfVHYI6HHDjEr%2BL7CqEtT5xQe7n4IBtkAEVHQHsQSgOibBldcbPNcy0AtvDasZ8wwUnFavnfDsA%2B2zwf8eeC5hr6mXaqXl1qbS6UYNho5RT%2BX%2BG9v8DMKR2ObtAPyAnMCuw%2BFPAoi%2BstvucmH8aphtdZ3ejVqEIlOxOvM%2BSwjfc5%2B3EuTGZqDrYt2h197VBuYqlLNpB6AUPziPuQIrtRItVzVHqZ2AdcKj0Q1XkBGhT2TRb%2BjiJKgYkdP3ItEwqDstXyvu0ITHjW1od1NxCxWlVI3vY55nDy7Nyt8xkwtEU%2BbWdCHEMIqLqL8GOqQBxegDiC2y%2Bp13DLt6fKpQmgO2a7V%2Fy5hxWlZeFK3iqVswg3byI1JAwHqE88HAmIXfE%2BkNMxog%2F23q6mCN1jVdKDULFrvSTH%2F8Od%2Bpu%2Fxu3rqynHz0XDHLTYTpS9khF0DYx1bz9sZ58N%2FNGV5L4ieQzW0URz4TbJ6p1C21VBAhqOCo0V8ocCJg63K7aQvB48wcd1vgOcv51qgc%2Fc2vgBZ364C9VFQ%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2012%3A00%3A04%20GMT&X-Amz-Signature=f40ac2bb25acf29abcb5b94a7aa9897cf2abb64054f75fb323a4c61967babf30'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276164863/OF6jgx1Lcp-camera-video-segment-driver-1743276167363.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOORQEFE7%2F20250331%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250331T040004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEDMaCXVzLXdlc3QtMiJIMEYCIQDYhldwn6xPibeo7UfO2hyaXchwZG%2BRUv3wHEvlAhzQ9wIhAKBCQOOXfTi9Kkihriivg46cyjPwSN9ea4juSW%2FfS77BKuYDCJz%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgwHR4J2QBwReOgfmJ8qugOPSbjSDvLSdXHFMzDd7rNCarOJd5s20MqO16y4tZIBDWKPK59sZrEyTMhwjPUPSLReD3kg2E1YPD%2Bvwc7OKmU7bmHhdFDUfIG7tKilS4GYE3JEvrsLY1pnH0X0VXlYEqm5ItI8QIHmvn29iLVQRQ1DWamIZa5%2BriHHXT9WhraXH2N4tTbrU8elcJX8DP%2FO8iLZCq5hMcBNdYb1O8WZTqxqCOZd1DjKZRtHgIQck5TdBZb0gwYLBPQzhamSl4kMMu%2FMH0vEeb%2B%2BZjz3r6f5V4I0Iy0gMIOfVHYI6HHDjEr%2BL7CqEtT5xQe7n4IBtkAEVHQHsQSgOibBldcbPNcy0AtvDasZ8wwUnFavnfDsA%2B2zwf8eeC5hr6mXaqXl1qbS6UYNho5RT%2BX%2BG9v8DMKR2ObtAPyAnMCuw%2BFPAoi%2BstvucmH8aphtdZ3ejVqEIlOxOvM%2BSwjfc5%2B3EuTGZqDrYt2h197VBuYqlLNpB6AUPziPuQIrtRItVzVHqZ2AdcKj0Q1XkBGhT2TRb%2BjiJKgYkdP3ItEwqDstXyvu0ITHjW1od1NxCxWlVI3vY55nDy7Nyt8xkwtEU%2BbWdCHEMIqLqL8GOqQBxegDiC2y%2Bp13DLt6fKpQmgO2a7V%2Fy5hxWlZeFK3iqVswg3byI1JAwHqE88HAmIXfE%2BkNMxog%2F23q6mCN1jVdKDULFrvSTH%2F8Od%2Bpu%2Fxu3rqynHz0XDHLTYTpS9khF0DYx1bz9sZ58N%2FNGV5L4ieQzW0URz4TbJ6p1C21VBAhqOCo0V8ocCJg63K7aQvB48wcd1vgOcv51qgc%2Fc2vgBZ364C9VFQ%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2012%3A00%3A04%20GMT&X-Amz-Signature=00ccd4b2e783d8252015d55a8675f53929e3b6440f068f23d6a673bb5e5ec4b0'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276108738/gAxSKLP636-camera-video-segment-driver-1743276111238.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOORQEFE7%2F20250331%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250331T040004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEDMaCXVzLXdlc3QtMiJIMEYCIQDYhldwn6xPibeo7UfO2hyaXchwZG%2BRUv3wHEvlAhzQ9wIhAKBCQOOXfTi9Kkihriivg46cyjPwSN9ea4juSW%2FfS77BKuYDCJz%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgwHR4J2QBwReOgfmJ8qugOPSbjSDvLSdXHFMzDd7rNCarOJd5s20MqO16y4tZIBDWKPK59sZrEyTMhwjPUPSLReD3kg2E1YPD%2Bvwc7OKmU7bmHhdFDUfIG7tKilS4GYE3JEvrsLY1pnH0X0VXlYEqm5ItI8QIHmvn29iLVQRQ1DWamIZa5%2BriHHXT9WhraXH2N4tTbrU8elcJX8DP%2FO8iLZCq5hMcBNdYb1O8WZTqxqCOZd1DjKZRtHgIQck5TdBZb0gwYLBPQzhamSl4kMMu%2FMH0vEeb%2B%2BZjz3r6f5V4I0Iy0gMIOfVHYI6HHDjEr%2BL7CqEtT5xQe7n4IBtkAEVHQHsQSgOibBldcbPNcy0AtvDasZ8wwUnFavnfDsA%2B2zwf8eeC5hr6mXaqXl1qbS6UYNho5RT%2BX%2BG9v8DMKR2ObtAPyAnMCuw%2BFPAoi%2BstvucmH8aphtdZ3ejVqEIlOxOvM%2BSwjfc5%2B3EuTGZqDrYt2h197VBuYqlLNpB6AUPziPuQIrtRItVzVHqZ2AdcKj0Q1XkBGhT2TRb%2BjiJKgYkdP3ItEwqDstXyvu0ITHjW1od1NxCxWlVI3vY55nDy7Nyt8xkwtEU%2BbWdCHEMIqLqL8GOqQBxegDiC2y%2Bp13DLt6fKpQmgO2a7V%2Fy5hxWlZeFK3iqVswg3byI1JAwHqE88HAmIXfE%2BkNMxog%2F23q6mCN1jVdKDULFrvSTH%2F8Od%2Bpu%2Fxu3rqynHz0XDHLTYTpS9khF0DYx1bz9sZ58N%2FNGV5L4ieQzW0URz4TbJ6p1C21VBAhqOCo0V8ocCJg63K7aQvB48wcd1vgOcv51qgc%2Fc2vgBZ364C9VFQ%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2012%3A00%3A04%20GMT&X-Amz-Signature=e52321d3b4570d26352661741eeef812dc3ce6ecd61b2e7c2ec1741bacbc24dd'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205341/1743265652740/hq4516Mzws-camera-video-segment-driver-1743265655240.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOORQEFE7%2F20250331%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250331T040004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEDMaCXVzLXdlc3QtMiJIMEYCIQDYhldwn6xPibeo7UfO2hyaXchwZG%2BRUv3wHEvlAhzQ9wIhAKBCQOOXfTi9Kkihriivg46cyjPwSN9ea4juSW%2FfS77BKuYDCJz%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgwHR4J2QBwReOgfmJ8qugOPSbjSDvLSdXHFMzDd7rNCarOJd5s20MqO16y4tZIBDWKPK59sZrEyTMhwjPUPSLReD3kg2E1YPD%2Bvwc7OKmU7bmHhdFDUfIG7tKilS4GYE3JEvrsLY1pnH0X0VXlYEqm5ItI8QIHmvn29iLVQRQ1DWamIZa5%2BriHHXT9WhraXH2N4tTbrU8elcJX8DP%2FO8iLZCq5hMcBNdYb1O8WZTqxqCOZd1DjKZRtHgIQck5TdBZb0gwYLBPQzhamSl4kMMu%2FMH0vEeb%2B%2BZjz3r6f5V4I0Iy0gMIOfVHYI6HHDjEr%2BL7CqEtT5xQe7n4IBtkAEVHQHsQSgOibBldcbPNcy0AtvDasZ8wwUnFavnfDsA%2B2zwf8eeC5hr6mXaqXl1qbS6UYNho5RT%2BX%2BG9v8DMKR2ObtAPyAnMCuw%2BFPAoi%2BstvucmH8aphtdZ3ejVqEIlOxOvM%2BSwjfc5%2B3EuTGZqDrYt2h197VBuYqlLNpB6AUPziPuQIrtRItVzVHqZ2AdcKj0Q1XkBGhT2TRb%2BjiJKgYkdP3ItEwqDstXyvu0ITHjW1od1NxCxWlVI3vY55nDy7Nyt8xkwtEU%2BbWdCHEMIqLqL8GOqQBxegDiC2y%2Bp13DLt6fKpQmgO2a7V%2Fy5hxWlZeFK3iqVswg3byI1JAwHqE88HAmIXfE%2BkNMxog%2F23q6mCN1jVdKDULFrvSTH%2F8Od%2Bpu%2Fxu3rqynHz0XDHLTYTpS9khF0DYx1bz9sZ58N%2FNGV5L4ieQzW0URz4TbJ6p1C21VBAhqOCo0V8ocCJg63K7aQvB48wcd1vgOcv51qgc%2Fc2vgBZ364C9VFQ%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2012%3A00%3A04%20GMT&X-Amz-Signature=3495b1f685214308e09b10b5b13ae6d87694171944e49d1f07fb9ac06038e98a'
$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743255999969/Yjq7ZuVo4p-camera-video-segment-driver-1743256002469.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOORQEFE7%2F20250331%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250331T040004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEDMaCXVzLXdlc3QtMiJIMEYCIQDYhldwn6xPibeo7UfO2hyaXchwZG%2BRUv3wHEvlAhzQ9wIhAKBCQOOXfTi9Kkihriivg46cyjPwSN9ea4juSW%2FfS77BKuYDCJz%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgwHR4J2QBwReOgfmJ8qugOPSbjSDvLSdXHFMzDd7rNCarOJd5s20MqO16y4tZIBDWKPK59sZrEyTMhwjPUPSLReD3kg2E1YPD%2Bvwc7OKmU7bmHhdFDUfIG7tKilS4GYE3JEvrsLY1pnH0X0VXlYEqm5ItI8QIHmvn29iLVQRQ1DWamIZa5%2BriHHXT9WhraXH2N4tTbrU8elcJX8DP%2FO8iLZCq5hMcBNdYb1O8WZTqxqCOZd1DjKZRtHgIQck5TdBZb0gwYLBPQzhamSl4kMMu%2FMH0vEeb%2B%2BZjz3r6f5V4I0Iy0gMIOfVHYI6HHDjEr%2BL7CqEtT5xQe7n4IBtkAEVHQHsQSgOibBldcbPNcy0AtvDasZ8wwUnFavnfDsA%2B2zwf8eeC5hr6mXaqXl1qbS6UYNho5RT%2BX%2BG9v8DMKR2ObtAPyAnMCuw%2BFPAoi%2BstvucmH8aphtdZ3ejVqEIlOxOvM%2BSwjfc5%2B3EuTGZqDrYt2h197VBuYqlLNpB6AUPziPuQIrtRItVzVHqZ2AdcKj0Q1XkBGhT2TRb%2BjiJKgYkdP3ItEwqDstXyvu0ITHjW1od1NxCxWlVI3vY55nDy7Nyt8xkwtEU%2BbWdCHEMIqLqL8GOqQBxegDiC2y%2Bp13DLt6fKpQmgO2a7V%2Fy5hxWlZeFK3iqVswg3byI1JAwHqE88HAmIXfE%2BkNMxog%2F23q6mCN1jVdKDULFrvSTH%2F8Od%2Bpu%2Fxu3rqynHz0XDHLTYTpS9khF0DYx1bz9sZ58N%2FNGV5L4ieQzW0URz4TbJ6p1C21VBAhqOCo0V8ocCJg63K7aQvB48wcd1vgOcv51qgc%2Fc2vgBZ364C9VFQ%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2012%3A00%3A04%20GMT&X-Amz-Signature=c7d1748c233a301bcacec4c6b55652c8b6fcc29518244e66fd7cf9cce84b981b'
